$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the UART peripheral row (row 5): rename UART6 -> UART3, update its
# description, and add a Notes entry describing the baud rate.
$ws.Range("A5").Value = "UART3"
$ws.Range("B5").Value = "Universal Asynchronous Recevier Transmitter 3"
$ws.Range("E5").Value = "Baud rate set to 115200 bps"

# Update the active selection to match the new edit location.
$ws.Range("F5").Select()
